$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = "2000 : 1"
$ws.Range("H2").Value = "500 : 1"
